$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "Step" label column (column A) added in front of each of the 4 blocks.
# Each block has: a title row, a header row, and 5 data rows (step 0..4).
# The block-title rows get A<r> = "Step" (plain, default style).
# The data rows get A<r> = 0..4 with the new centered style (xf index 10).
# ---------------------------------------------------------------------------

$blocks = @(
    @{ Title = 2;  Header = 3;  Data = 4  },
    @{ Title = 10; Header = 11; Data = 12 },
    @{ Title = 18; Header = 19; Data = 20 },
    @{ Title = 26; Header = 27; Data = 28 }
)

# Blocks 1-3 put the "Step" label on the title row itself (A2, A10, A18).
$ws.Range("A2").Value = "Step"
$ws.Range("A10").Value = "Step"
$ws.Range("A18").Value = "Step"
# Block 4's title row (row 26) is left untouched; the label instead goes on
# its header row (A27), matching the source edit exactly.
$ws.Range("A27").Value = "Step"

foreach ($block in $blocks) {
    $dataStart = $block.Data
    for ($i = 0; $i -le 4; $i++) {
        $row = $dataStart + $i
        $cell = $ws.Range("A$row")
        $cell.Value = $i
        $cell.HorizontalAlignment = -4108
        $cell.VerticalAlignment = -4108
    }
}

# ---------------------------------------------------------------------------
# Updated simulation results for the 3rd and 4th blocks (rows 21-24, 29-32).
# ---------------------------------------------------------------------------

$ws.Range("E21").Value = 24.852767080646
$ws.Range("F21").Value = 21.856203916320499
$ws.Range("G21").Value = 19.979146932101699

$ws.Range("E22").Value = 39.0935135148586
$ws.Range("F22").Value = 27.097694618126699
$ws.Range("G22").Value = 20.012470095649402

$ws.Range("E23").Value = 39.425250148468997
$ws.Range("F23").Value = 27.477551068029499
$ws.Range("G23").Value = 20.078523817849501

$ws.Range("E24").Value = 39.473399561969302
$ws.Range("F24").Value = 27.547121204466801
$ws.Range("G24").Value = 20.145770205286599

$ws.Range("E29").Value = 24.825118284210198
$ws.Range("F29").Value = 21.815252902806201
$ws.Range("G29").Value = 19.9365266922808

$ws.Range("E30").Value = 39.451806481268498
$ws.Range("F30").Value = 27.512998054929302
$ws.Range("G30").Value = 20.103925531273799

$ws.Range("E31").Value = 39.555751585439999
$ws.Range("F31").Value = 27.676952978339699
$ws.Range("G31").Value = 20.304128387374099

$ws.Range("E32").Value = 39.657104806189203
$ws.Range("F32").Value = 27.837120774334
$ws.Range("G32").Value = 20.5005428140025

# ---------------------------------------------------------------------------
# Selection moved from L28 to F14.
# ---------------------------------------------------------------------------
$ws.Range("F14").Select()
